# feat: modify loginId type
#
# - Column A header: "내선번호" (extension number) -> "로그인아이디" (login ID)
# - A2/A3: numeric extension numbers -> text login IDs (loginId1 / loginId2),
#   right aligned like the old hyperlink-styled header cells
# - D2/D3: drop the (stale/incorrect) mailto hyperlinks and their special
#   hyperlink styling, keep the underlying email text untouched
# - E2/E3: reformat the sample phone number with dashes

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: "내선번호" -> "로그인아이디" -------------------------------
$ws.Range("A1").Value = "로그인아이디"

# A2/A3 used to hold plain numeric extension numbers; they now hold text
# login ids, right aligned (matching the look of the other right-aligned
# "하이퍼링크" style cells used elsewhere in the template).
$ws.Range("A2").Value = "loginId1"
$ws.Range("A3").Value = "loginId2"
$ws.Range("A2:A3").HorizontalAlignment = -4152

# --- D2/D3: remove the stray hyperlinks + their styling --------------------
$ws.Hyperlinks.Delete()

# Reset D2/D3 back to the workbook's default (unstyled) formatting, the same
# way the other plain-text cells (e.g. F3) in the sheet look, without
# introducing a brand-new style definition.
$ws.Range("F3").Copy()
$ws.Range("D2:D3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- E2/E3: normalize the sample phone number format -----------------------
$ws.Range("E2").Value = "010-9999-9999"
$ws.Range("E3").Value = "010-9999-9999"

# --- Misc: restore the selection state left behind in the saved file -------
$excel.Goto($ws.Range("D4:D5"))
